# Changes regarding Alex & Ansgar
# Add an "Accession ID" column in front of the existing data, drop the
# Accuracy / Relevance / Coherence columns, and fix the AGILYSYS INC name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before column A so Company_Name -> B, Year -> C.
$ws.Columns.Item(1).Insert()

# 2) New header + Accession ID values (row order matches the existing rows).
$ws.Range("A1").Value = "Accession ID"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A2").Value = "0000897101-17-000402"
$ws.Range("A3").Value = "0001615774-19-005785"
$ws.Range("A4").Value = "0001326380-18-000033"
$ws.Range("A5").Value = "0001683168-21-001026"
$ws.Range("A6").Value = "0000768835-18-000030"
$ws.Range("A7").Value = "0001309108-23-000022"
$ws.Range("A8").Value = "0000900075-23-000034"
$ws.Range("A9").Value = "0001731122-20-000381"
$ws.Range("A10").Value = "0000078749-17-000021"
$ws.Range("A11").Value = "0001509991-21-000031"

# 3) Fix the AGILYSYS INC name (drop the stray leading tab) - now in column B.
$ws.Range("B10").Value = "AGILYSYS INC"

# 4) Drop the old Accuracy / Relevance / Coherence columns (now D:F).
$ws.Range("D1:F11").Delete()

# 5) Re-apply sane column widths for the new A:C layout.
$ws.Columns.Item(1).ColumnWidth = 20.17
$ws.Columns.Item(2).ColumnWidth = 30
$ws.Columns.Item(3).ColumnWidth = 17.6

$ws.Range("B22").Select()
